# Sync attendance_reports: rotate the "Recorded By" (column G) author list
# on the "Session Analysis Results" sheet so the last-listed recorder moves
# to the front of the comma-separated list, for the specific rows touched
# upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G ("Recorded By") whose value needs to change.
$rowsToUpdate = @(
    2,3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,32,36,37,38,39,40,
    41,43,44,45,46,47,48,50,52,54,55,58,62,63,64,65,66,67,69,70,71,72,73,74,
    76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,
    119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153
)

foreach ($row in $rowsToUpdate) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $current = [string]$cell.Text

    # Split on ", " to preserve the original separator formatting.
    $parts = $current -split ", "

    if ($parts.Count -gt 1) {
        # Move the last author in the list to the front.
        $last = $parts[$parts.Count - 1]
        $rest = $parts[0..($parts.Count - 2)]
        $newParts = @($last) + @($rest)
        $newValue = [string]::Join(", ", $newParts)
        $cell.Value = $newValue
    }
}
